$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark that currently sits after "MP73010" in the title line.
#    We'll re-add an equivalent bookmark later in the ">>> ... >>>" paragraph.
foreach ($bm in @($d.Bookmarks)) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# 2. Rewrite the ">>  > your stuff after this line >>>" paragraph into the new text,
#    inserting the _GoBack bookmark in its new location.
$d.Content.Find.Execute(">>>  your stuff after this line >>>", $true, $false, $false, $false, $false, $true, 1, $false, ">>>  my first fork change done. thanks<<<BOOKMARK>>> >>>", 2)

# Find the paragraph range containing our placeholder marker and split it up, inserting
# the bookmark where <<<BOOKMARK>>> currently is.
$rng = $d.Content
$found = $rng.Find.Execute("<<<BOOKMARK>>>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $d.Bookmarks.Add("_GoBack", $rng)
    $rng.Text = ""
}
